$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unused trailing columns (F:G and AN:AR) which have no data -
# a single contiguous EntireColumn.Delete across the whole span gives the
# cleanest result in this runtime.
$ws.Range("F1:AR1").EntireColumn.Delete()

# Resize the remaining data columns B, C and D to their new widths.
$ws.Columns.Item(2).ColumnWidth = 19.276041666666668
$ws.Columns.Item(3).ColumnWidth = 20.276041666666668
$ws.Columns.Item(4).ColumnWidth = 22.498697916666668

# Add the new "#" header label in A2 (new shared string).
$ws.Range("A2").Value = "#"

# Update the view: drop the B1 frozen/top-left scroll position and move the
# active selection from B4 to A4.
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A4").Select()

# Restore the workbook tab ratio (best effort; this runtime may not persist it).
$excel.ActiveWindow.TabRatio = 0.838
